$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: add a new "Update doc gen tools and improve doc gen process." list
# item right after the "Compile-fail tests for a Call." bullet, and relocate
# the "_GoBack" bookmark (Word's "last edit location" marker) into the middle
# of the new sentence, matching the point where the author's cursor ended up.
# ---------------------------------------------------------------------------

$sourceIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Compile-fail tests for a Call.") {
        $sourceIdx = $i
        break
    }
}

if ($sourceIdx -gt 0) {
    $srcPara = $d.Paragraphs($sourceIdx)
    $srcRange = $srcPara.Range
    $srcRange.Collapse(0)
    $srcRange.InsertParagraphAfter()

    $newIdx = $sourceIdx + 1

    $newRange = $d.Paragraphs($newIdx).Range
    $newRange.InsertAfter("Update doc gen tools and improve doc gen process.")

    # Put the bookmark back where the author's cursor last was: right before
    # "process." (i.e. after "...improve doc gen ").
    $locateRange = $d.Paragraphs($newIdx).Range
    $locateRange.Find.Execute("process.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $bmPoint = $d.Range($locateRange.Start, $locateRange.Start)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}

# ---------------------------------------------------------------------------
# Hunks 2-5: the lastRenderedPageBreak marker shifts one bullet earlier in
# two different lists (a side effect of the content added above pushing the
# page break back by a couple of lines). Move the (empty, textless) marker
# from the later paragraph to the earlier one in each pair, leaving
# everything else about those paragraphs untouched.
# ---------------------------------------------------------------------------

function Move-LastRenderedPageBreak($fromText, $toText) {
    $fromIdx = 0
    $toIdx = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($fromIdx -eq 0 -and $t -like $fromText) { $fromIdx = $i }
        if ($toIdx -eq 0 -and $t -like $toText) { $toIdx = $i }
    }
    if ($fromIdx -eq 0 -or $toIdx -eq 0) {
        return
    }

    $fromPara = $d.Paragraphs($fromIdx)
    $fromRuns = $fromPara.Range
    # Remove it from the "from" paragraph by rebuilding that run without it.
    # (Read its current formatted text back out so nothing else changes.)

    $toPara = $d.Paragraphs($toIdx)

    # Use a tiny, local InsertXML replace scoped exactly to each paragraph's
    # own Range so neighboring paragraphs are left completely alone.
    return @($fromPara, $toPara)
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Pair 1: "Fastcall function hooking" (gains it) / "VEH hooking" (loses it) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Fastcall function hooking*") {
        $full = $d.Paragraphs($i).Range
        $xml = '<w:p ' + $wNs + ' w:rsidR="00BA3F8E" w:rsidRDefault="00BA3F8E" w:rsidP="00BA3F8E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Fastcall function hooking (ecx, edx preservation). </w:t></w:r></w:p>'
        $full.InsertXML($xml)
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*VEH hooking*") {
        $full = $d.Paragraphs($i).Range
        $xml = '<w:p ' + $wNs + ' w:rsidR="00312CFE" w:rsidRPr="008B06FC" w:rsidRDefault="00312CFE" w:rsidP="00312CFE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:t>VEH hooking (both INT3 and DR).</w:t></w:r></w:p>'
        $full.InsertXML($xml)
        break
    }
}

# --- Pair 2: "CLR runtime directory support." (gains it) / "Full support for
# writing back to PE file..." (loses it) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*CLR runtime directory support*") {
        $full = $d.Paragraphs($i).Range
        $xml = '<w:p ' + $wNs + ' w:rsidR="009420FA" w:rsidRPr="008B06FC" w:rsidRDefault="009420FA" w:rsidP="009420FA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:lastRenderedPageBreak/><w:t>CLR runtime directory support.</w:t></w:r></w:p>'
        $full.InsertXML($xml)
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Full support for writing back to PE file*") {
        $full = $d.Paragraphs($i).Range
        $xml = '<w:p ' + $wNs + ' w:rsidR="009420FA" w:rsidRPr="008B06FC" w:rsidRDefault="009420FA" w:rsidP="008F1E4C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:t>Full support for writing back to PE file, including automatically performing adjustments where required to fit in new data or remove unnecessary space.</w:t></w:r></w:p>'
        $full.InsertXML($xml)
        break
    }
}

Write-Output "Done."
